$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C6").Value = "222222222"
$ws.Range("I6").Value = "/201761856"
$ws.Range("A12").Value = "2017-10-07"
$ws.Range("F12").Value = "212"
$ws.Range("H12").Value = "12"
$ws.Range("J12").Value = "Femenino"
$ws.Range("A18").Value = "1111"
$ws.Range("F18").Value = "1111"
$ws.Range("A24").Value = "23/10/2017"
$ws.Range("C24").Value = "10:42:11"
